$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (ALC)
$ws.Range("H17").Value = 1377.4474
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1377.4474
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4132.3422
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -4468.3422

# Row 62 (ALC)
$ws.Range("H62").Value = 2032.5
$ws.Range("I62").Value = 1839
$ws.Range("K62").Value = 1839
$ws.Range("M62").Value = -1215

# Row 64 (ALC)
$ws.Range("H64").Value = 3150
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# Row 65 (ALC)
$ws.Range("H65").Value = 2032.5
$ws.Range("I65").Value = 1839
$ws.Range("K65").Value = 9195
$ws.Range("M65").Value = -6075

# Row 67 (ALC)
$ws.Range("H67").Value = 3150
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# Row 103 (ALC)
$ws.Range("H103").Value = 999.3913
$ws.Range("J103").Value = 998
$ws.Range("L103").Value = 2994
$ws.Range("N103").Value = -4166

# Row 132 (ALC)
$ws.Range("H132").Value = 1125
$ws.Range("I132").Value = 938
$ws.Range("K132").Value = 2814
$ws.Range("M132").Value = -284

# Row 138 (ALC)
$ws.Range("H138").Value = 2267.1
$ws.Range("I138").Value = 2196.9355
$ws.Range("K138").Value = 6590.806500000001
$ws.Range("M138").Value = -1450.806500000001

# Row 141 (ALC)
$ws.Range("H141").Value = 3941.45
$ws.Range("I141").Value = 3188
$ws.Range("K141").Value = 9564
$ws.Range("M141").Value = -4384

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Range("H61").Value = 3582.9312
$ws.Range("I61").Value = 2658.8667
$ws.Range("J61").Value = 4573
$ws.Range("K61").Value = 2658.8667
$ws.Range("L61").Value = 4573
$ws.Range("M61").Value = -2446.8667
$ws.Range("N61").Value = -4997

# Row 88 (ARM)
$ws.Range("H88").Value = 4381.364
$ws.Range("I88").Value = 3199
$ws.Range("J88").Value = 5057
$ws.Range("K88").Value = 3199
$ws.Range("L88").Value = 5057
$ws.Range("M88").Value = -2793
$ws.Range("N88").Value = -5869

# Row 91 (ARM)
$ws.Range("H91").Value = 4381.364
$ws.Range("I91").Value = 3199
$ws.Range("J91").Value = 5057
$ws.Range("K91").Value = 3199
$ws.Range("L91").Value = 5057
$ws.Range("M91").Value = -1795
$ws.Range("N91").Value = -7865

# Row 136 (ARM)
$ws.Range("H136").Value = 3582.9312
$ws.Range("I136").Value = 2658.8667
$ws.Range("J136").Value = 4573
$ws.Range("K136").Value = 7976.6001
$ws.Range("L136").Value = 13719
$ws.Range("M136").Value = -5426.6001
$ws.Range("N136").Value = -18819

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 2998.4119
$ws.Range("I31").Value = 5749
$ws.Range("J31").Value = 2631.6667
$ws.Range("K31").Value = 5749
$ws.Range("L31").Value = 2631.6667
$ws.Range("M31").Value = -5454
$ws.Range("N31").Value = -3221.6667

# Row 34 (CRP)
$ws.Range("H34").Value = 2998.4119
$ws.Range("I34").Value = 5749
$ws.Range("J34").Value = 2631.6667
$ws.Range("K34").Value = 5749
$ws.Range("L34").Value = 2631.6667
$ws.Range("M34").Value = -5547
$ws.Range("N34").Value = -3035.6667

# Row 94 (CRP)
$ws.Range("H94").Value = 1007.8
$ws.Range("I94").Value = 955.5
$ws.Range("J94").Value = 1020.875
$ws.Range("K94").Value = 955.5
$ws.Range("L94").Value = 1020.875
$ws.Range("M94").Value = -504.5
$ws.Range("N94").Value = -1922.875

# Row 132 (CRP)
$ws.Range("H132").Value = 1895.591
$ws.Range("I132").Value = 1184.8
$ws.Range("J132").Value = 3418.7144
$ws.Range("K132").Value = 3554.4
$ws.Range("L132").Value = 10256.1432
$ws.Range("M132").Value = -1024.4
$ws.Range("N132").Value = -15316.1432

$ws = $wb.Worksheets.Item("CUL")
# Row 107 (CUL)
$ws.Range("H107").Value = 901.5454999999999
$ws.Range("I107").Value = 398
$ws.Range("J107").Value = 925.5238000000001
$ws.Range("K107").Value = 1194
$ws.Range("L107").Value = 2776.5714
$ws.Range("M107").Value = 726
$ws.Range("N107").Value = -6616.571400000001

# Row 139 (CUL)
$ws.Range("H139").Value = 12046.5
$ws.Range("I139").Value = 15937.857
$ws.Range("K139").Value = 47813.571
$ws.Range("M139").Value = -42673.571

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Range("H70").Value = 4931.25
$ws.Range("I70").Value = 4737.5
$ws.Range("K70").Value = 4737.5
$ws.Range("M70").Value = -4467.5

# Row 73 (GSM)
$ws.Range("H73").Value = 4931.25
$ws.Range("I73").Value = 4737.5
$ws.Range("K73").Value = 4737.5
$ws.Range("M73").Value = -3801.5

# Row 98 (GSM)
$ws.Range("H98").Value = 19007.5
$ws.Range("J98").Value = 19007.5
$ws.Range("L98").Value = 19007.5
$ws.Range("N98").Value = -24997.5

# Row 102 (GSM)
$ws.Range("H102").Value = 2982.077
$ws.Range("I102").Value = 2980.5833
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2980.5833
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -1358.5833
$ws.Range("N102").Value = -6244

# Row 122 (GSM)
$ws.Range("H122").Value = 1765.7222
$ws.Range("I122").Value = 1429.6154
$ws.Range("K122").Value = 4288.8462
$ws.Range("M122").Value = -1838.8462

# Row 132 (GSM)
$ws.Range("H132").Value = 3207520.5
$ws.Range("I132").Value = 4809357.5
$ws.Range("J132").Value = 3846.5
$ws.Range("K132").Value = 14428072.5
$ws.Range("L132").Value = 11539.5
$ws.Range("M132").Value = -14425542.5
$ws.Range("N132").Value = -16599.5

# Row 139 (GSM)
$ws.Range("H139").Value = 55227.1
$ws.Range("J139").Value = 55227.1
$ws.Range("L139").Value = 55227.1
$ws.Range("N139").Value = -65507.1

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("H16").Value = 3273.5293
$ws.Range("I16").Value = 3786.3333
$ws.Range("J16").Value = 2696.625
$ws.Range("K16").Value = 3786.3333
$ws.Range("L16").Value = 2696.625
$ws.Range("M16").Value = -3616.3333
$ws.Range("N16").Value = -3036.625

# Row 22 (LTW)
$ws.Range("H22").Value = 2853.7
$ws.Range("J22").Value = 1864.5
$ws.Range("L22").Value = 1864.5
$ws.Range("N22").Value = -2454.5

# Row 27 (LTW)
$ws.Range("H27").Value = 2853.7
$ws.Range("J27").Value = 1864.5
$ws.Range("L27").Value = 1864.5
$ws.Range("N27").Value = -2078.5

# Row 46 (LTW)
$ws.Range("H46").Value = 2597.8333
$ws.Range("J46").Value = 3224.8572
$ws.Range("L46").Value = 3224.8572
$ws.Range("N46").Value = -3600.8572

# Row 55 (LTW)
$ws.Range("H55").Value = 256.66666
$ws.Range("I55").Value = 195.91667
$ws.Range("J55").Value = 499.66666
$ws.Range("K55").Value = 195.91667
$ws.Range("L55").Value = 499.66666
$ws.Range("M55").Value = -22.91667000000001
$ws.Range("N55").Value = -845.66666

# Row 136 (LTW)
$ws.Range("H136").Value = 2923.1
$ws.Range("J136").Value = 4443.5
$ws.Range("L136").Value = 13330.5
$ws.Range("N136").Value = -18430.5

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (WVR)
$ws.Range("H96").Value = 9455.444
$ws.Range("I96").Value = 600
$ws.Range("K96").Value = 600
$ws.Range("M96").Value = 773

# Row 100 (WVR)
$ws.Range("H100").Value = 614.375
$ws.Range("I100").Value = 435.83334
$ws.Range("K100").Value = 871.66668
$ws.Range("M100").Value = -330.66668

# Row 112 (WVR)
$ws.Range("H112").Value = 14001
$ws.Range("J112").Value = 14001
$ws.Range("L112").Value = 14001
$ws.Range("N112").Value = -16955
